$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the affected columns so Excel does not
# auto-convert the numeric-looking / percent-looking strings into
# actual numbers when the new values are assigned.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "278.69"
$ws.Range("E2").Value = "6.65%"
$ws.Range("G2").Value = "2"

$ws.Range("D3").Value = "27.21"
$ws.Range("E3").Value = "-2.28%"
$ws.Range("G3").Value = "2"

$ws.Range("D4").Value = "4.787"
$ws.Range("E4").Value = "1.96%"
$ws.Range("G4").Value = "2"

$ws.Range("D5").Value = "0.06252"
$ws.Range("E5").Value = "0.42%"
$ws.Range("G5").Value = "2"

$ws.Range("D6").Value = "6.810"
$ws.Range("E6").Value = "1.04%"
$ws.Range("G6").Value = "2"

$ws.Range("D7").Value = "0.8711"
$ws.Range("E7").Value = "2.47%"
$ws.Range("G7").Value = "2"

$ws.Range("D8").Value = "0.9505"
$ws.Range("E8").Value = "4.52%"
$ws.Range("G8").Value = "2"

$ws.Range("D9").Value = "0.1460"
$ws.Range("E9").Value = "4.28%"
$ws.Range("G9").Value = "2"

$ws.Range("D10").Value = "0.05260"
$ws.Range("E10").Value = "10.09%"
$ws.Range("G10").Value = "2"

$ws.Range("D11").Value = "0.07239"
$ws.Range("E11").Value = "2.06%"
$ws.Range("G11").Value = "2"

$ws.Range("D12").Value = "0.03124"
$ws.Range("E12").Value = "0.17%"
$ws.Range("G12").Value = "2"

$ws.Range("D13").Value = "0.09047"
$ws.Range("E13").Value = "-0.04%"
$ws.Range("G13").Value = "2"

$ws.Range("D14").Value = "0.001546"
$ws.Range("E14").Value = "0.76%"
$ws.Range("G14").Value = "2"

$ws.Range("D15").Value = "0.0006267"
$ws.Range("E15").Value = "1.74%"
$ws.Range("G15").Value = "2"

$ws.Range("D16").Value = "0.005966"
$ws.Range("E16").Value = "-0.87%"
$ws.Range("G16").Value = "2"

$ws.Range("D17").Value = "3.469"
$ws.Range("E17").Value = "0.17%"
$ws.Range("G17").Value = "2"

$ws.Range("D18").Value = "3.252"
$ws.Range("E18").Value = "2.48%"
$ws.Range("G18").Value = "2"

$ws.Range("G19").Value = "2"

$ws.Range("E20").Value = "-0.64%"
$ws.Range("G20").Value = "2"

$ws.Range("E21").Value = "-0.31%"
$ws.Range("G21").Value = "2"

$ws.Range("D22").Value = "3.827"
$ws.Range("E22").Value = "-6.79%"
$ws.Range("G22").Value = "2"

$ws.Range("D23").Value = "0.04293"
$ws.Range("E23").Value = "1.09%"
$ws.Range("G23").Value = "2"

$ws.Range("E24").Value = "-3.84%"
$ws.Range("G24").Value = "2"

$ws.Range("D25").Value = "0.004225"
$ws.Range("E25").Value = "3.29%"
$ws.Range("G25").Value = "2"

$ws.Range("D26").Value = "0.0001196"
$ws.Range("E26").Value = "-0.46%"
$ws.Range("G26").Value = "2"

$ws.Range("E27").Value = "18.15%"
$ws.Range("G27").Value = "2"

$ws.Range("G28").Value = "2"

$ws.Range("G29").Value = "2"

$ws.Range("G30").Value = "2"

$ws.Range("G31").Value = "2"

$ws.Range("G32").Value = "2"

$ws.Range("G33").Value = "2"

$ws.Range("G34").Value = "2"

$ws.Range("G35").Value = "2"

$ws.Range("G36").Value = "2"

$ws.Range("G37").Value = "2"

$ws.Range("G38").Value = "2"

$ws.Range("G39").Value = "2"

$ws.Range("D40").Value = "0.04040"
$ws.Range("E40").Value = "4.21%"
$ws.Range("G40").Value = "2"

$ws.Range("D41").Value = "0.006188"
$ws.Range("E41").Value = "49.92%"
$ws.Range("G41").Value = "2"

$ws.Range("D42").Value = "0.1140"
$ws.Range("E42").Value = "2.59%"
$ws.Range("G42").Value = "2"

$ws.Range("D43").Value = "0.002117"
$ws.Range("E43").Value = "-3.86%"
$ws.Range("G43").Value = "2"

$ws.Range("D44").Value = "0.01276"
$ws.Range("E44").Value = "-4.29%"
$ws.Range("G44").Value = "2"

$ws.Range("D45").Value = "0.00005120"
$ws.Range("E45").Value = "-0.72%"
$ws.Range("G45").Value = "2"

$ws.Range("D46").Value = "0.00000000747"
$ws.Range("E46").Value = "-0.45%"
$ws.Range("G46").Value = "2"

$ws.Range("D47").Value = "1.494"
$ws.Range("E47").Value = "2,671.67%"
$ws.Range("G47").Value = "2"

$ws.Range("D48").Value = "0.02981"
$ws.Range("E48").Value = "-12.43%"
$ws.Range("G48").Value = "2"

$ws.Range("E49").Value = "-0.45%"
$ws.Range("G49").Value = "2"

$ws.Range("E50").Value = "-0.45%"
$ws.Range("G50").Value = "2"

$ws.Range("G51").Value = "2"
